# Docker image standalone successful
#
# This script reproduces the authored edit to testdata.xlsx:
#   1. The "test_suite" sheet gains a "ProductPage" test-case row (inserted
#      before the existing "loginTest" row) and a "productPage" test-case
#      row appended after it. Both new A-column cells use the same
#      non-underlined black Arial font that the existing "LoginTest" link
#      cell's underline style was derived from.
#   2. A new, empty "Sheet3" worksheet is appended after "test_suite".
#   3. Selection/active-cell bookkeeping is restored so "test_suite" stays
#      the active tab (as it was before the edit) while the new sheet gets
#      its own default selection.

$wb = $excel.ActiveWorkbook

$wsSuite = $wb.Worksheets.Item(2)

# --- test_suite: push the existing "loginTest" case down to row 4 ---
$wsSuite.Range("A4").Value = "loginTest"
$wsSuite.Range("B4").Value = "Y"

# --- test_suite: row 3 becomes the new "ProductPage" case ---
$wsSuite.Range("A3").Value = "ProductPage"
$wsSuite.Range("B3").Value = "Y"
$wsSuite.Range("A3").Font.Name = "Arial"
$wsSuite.Range("A3").Font.Size = 10
$wsSuite.Range("A3").Font.Color = 0x000000
$wsSuite.Range("A3").Font.Underline = $false

# --- test_suite: row 5 is the new "productPage" case ---
$wsSuite.Range("A5").Value = "productPage"
$wsSuite.Range("B5").Value = "Y"
$wsSuite.Range("A5").Font.Name = "Arial"
$wsSuite.Range("A5").Font.Size = 10
$wsSuite.Range("A5").Font.Color = 0x000000
$wsSuite.Range("A5").Font.Underline = $false

# --- add the new, empty "Sheet3" after "test_suite" ---
$wsNew = $wb.Worksheets.Add($null, $wsSuite)
$wsNew.Name = "Sheet3"
$wsNew.Range("A2").Select() | Out-Null

# --- restore test_suite as the active sheet/selection ---
$wsSuite.Activate() | Out-Null
$wsSuite.Range("C14").Select() | Out-Null
